$wb = $excel.ActiveWorkbook

# --- Sheet: Neodymium ---
$ws = $wb.Worksheets.Item("Neodymium")

$ws.Range("C2").Value = [double]"4.05721109302746E-09"
$ws.Range("D2").Value = [double]"0.001664791307295803"
$ws.Range("E2").Value = [double]"0.004065106765944203"

$ws.Range("B3").Value = [double]"4.691044125953377E-15"
$ws.Range("C3").Value = [double]"1.512131801220274E-05"
$ws.Range("D3").Value = [double]"0.001058606147983076"
$ws.Range("E3").Value = [double]"0.003596305710095619"

$ws.Range("B4").Value = [double]"7.321669688613381E-17"
$ws.Range("C4").Value = [double]"3.740244660593755E-06"
$ws.Range("D4").Value = [double]"0.0007931982937362548"
$ws.Range("E4").Value = [double]"0.003171595483280042"

$ws.Range("C5").Value = [double]"3.409138439013245E-12"
$ws.Range("D5").Value = [double]"1.223194692503425E-05"
$ws.Range("E5").Value = [double]"0.0002413018472889734"

# --- Sheet: Copper ---
$ws = $wb.Worksheets.Item("Copper")

$ws.Range("B2").Value = [double]"3.013444709329602E-07"
$ws.Range("C2").Value = [double]"0.0001393251862087163"
$ws.Range("D2").Value = [double]"0.1266827953441403"
$ws.Range("E2").Value = [double]"0.3498475259112285"

$ws.Range("B3").Value = [double]"2.048071879604898E-06"
$ws.Range("C3").Value = [double]"0.002762933903811979"
$ws.Range("D3").Value = [double]"0.075856749380545"
$ws.Range("E3").Value = [double]"0.2421667282385033"

$ws.Range("B4").Value = [double]"6.08104089514155E-06"
$ws.Range("C4").Value = [double]"0.0003599297232404902"
$ws.Range("D4").Value = [double]"0.04738308718200741"
$ws.Range("E4").Value = [double]"0.2328833888754995"

$ws.Range("B5").Value = [double]"1.909822742276567E-06"
$ws.Range("C5").Value = [double]"0.0008394655527376822"
$ws.Range("D5").Value = [double]"0.08542358177263763"
$ws.Range("E5").Value = [double]"0.2358540537555117"

# --- Sheet: Raw silicon ---
$ws = $wb.Worksheets.Item("Raw silicon")

$ws.Range("B2").Value = [double]"5.031589966768446E-08"
$ws.Range("C2").Value = [double]"3.79581360833385E-06"
$ws.Range("D2").Value = [double]"0.002641707319413292"
$ws.Range("E2").Value = [double]"0.01383839401004388"

$ws.Range("B3").Value = [double]"5.367099139082082E-08"
$ws.Range("C3").Value = [double]"1.593037759322987E-05"
$ws.Range("D3").Value = [double]"0.00152498881473661"
$ws.Range("E3").Value = [double]"0.006174302220107622"

$ws.Range("B4").Value = [double]"3.443534520072915E-07"
$ws.Range("C4").Value = [double]"4.120808690356009E-06"
$ws.Range("D4").Value = [double]"0.001111119351128833"
$ws.Range("E4").Value = [double]"0.006154403591514441"

$ws.Range("B5").Value = [double]"1.848515343462404E-07"
$ws.Range("C5").Value = [double]"5.196111533190502E-06"
$ws.Range("D5").Value = [double]"0.002361937477439513"
$ws.Range("E5").Value = [double]"0.00981617050798976"
